$d = $word.ActiveDocument

function Set-RunXml($findText, $matchCase, $innerXml) {
    $found = $d.Content
    $found.Find.Execute($findText, $matchCase, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    # Re-seat a fresh Range over the found span - reusing the Find-owning Range
    # object directly with InsertXML leaves stray leading characters behind.
    $rng = $d.Range($found.Start, $found.End)
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# 1. Welcome -> Добро пожаловать  (needs xml:space="preserve" added)
$run1 = '<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:eastAsia="Times New Roman"/><w:b/><w:color w:val="2D2D2D"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Добро пожаловать</w:t></w:r>'
Set-RunXml "Welcome" $false $run1

# 2. SmartCash is pursuing ... -> Russian translation (trailing double space keeps preserve automatically)
$old2 = "SmartCash is pursuing a very ambitious development schedule and is constantly improving. Do you have a talent that could help? Join us on our community Discord and find ways to contribute."
$new2 = "SmartCash следует принципам постоянного развития и строгого выполнения намеченных планов. Если у вас есть талант, который может быть востребован, присоединяйтесь к нашему сообществу в Discord и узнайте, чем вы можете быть полезны.  "
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# 3. WHAT IS SMARTCASH? -> Что такое SmartCash?  (needs xml:space="preserve" retained)
$run3 = '<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:eastAsia="Times New Roman"/><w:caps/><w:color w:val="F4B517"/><w:spacing w:val="15"/><w:kern w:val="36"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Что такое SmartCash?</w:t></w:r>'
Set-RunXml "WHAT IS SMARTCASH?" $false $run3

# 4. SmartCash is a community governance... -> Russian translation with embedded newline (needs xml:space="preserve")
$old4 = "SmartCash is a community governance, cooperation & growth focused blockchain based currency & a decentralized economy."
$new4 = "SmartCash – это не просто криптовалюта на основе технологии блокчейн, но децентрализованная экономическая система,`nуправляемая сообществом для взаимовыгодного сотрудничества и роста."
$run4 = '<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">' + $new4 + '</w:t></w:r>'
Set-RunXml $old4 $false $run4

# 5. SMARTHIVE -> SmartHive (exact case match needed; no preserve either side)
$d.Content.Find.Execute("SMARTHIVE", $true, $false, $false, $false, $false, $true, 1, $false, "SmartHive", 2) | Out-Null
